$wb = $excel.ActiveWorkbook

# --- Append newly logged play-by-play yardage/return data (Week 15 log + Week 16 sim) ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value2 = $ws.Range("B2").Value2 + " 2 1 10 3 0 7 2 7 0 0 1 2 3 -3 6 -5"

$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B3").Value2 = $ws.Range("B3").Value2 + " 16 4 6 3 -6 11 14 7 19 5 5 17 12 7 5 8 4"

$ws = $wb.Worksheets.Item("YDS")
$ws.Range("C2").Value2 = $ws.Range("C2").Value2 + " 3 6 2 1 3 8 9 2 1 2 3 6 0 3 9 1 2 0 9 12 1 6 3 10 7 2 4 4 15 4 5 2 11 0 20 6 7 3 1 4 3"

$ws = $wb.Worksheets.Item("YDS")
$ws.Range("C3").Value2 = $ws.Range("C3").Value2 + " 7 5 11 12 1 4 8 5 18 7 9 13 11 5 -4 -3 3 5 12 5 9 8"

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B4").Value2 = $ws.Range("B4").Value2 + " 64 59"

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B5").Value2 = $ws.Range("B5").Value2 + " 21 19"

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B6").Value2 = $ws.Range("B6").Value2 + " 20 18"

$ws = $wb.Worksheets.Item("ST")
$ws.Range("D3").Value2 = $ws.Range("D3").Value2 + " 27 51 32 40 49"

$ws = $wb.Worksheets.Item("ST")
$ws.Range("D4").Value2 = $ws.Range("D4").Value2 + " 0 55 0 0 0"

$ws = $wb.Worksheets.Item("ST")
$ws.Range("D5").Value2 = $ws.Range("D5").Value2 + " 17 11 11"

# --- Update aggregate numeric totals on each stats sheet ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("B2").Value2 = 11
$ws.Range("C2").Value2 = 318
$ws.Range("D2").Value2 = 27
$ws.Range("G2").Value2 = 85
$ws.Range("J2").Value2 = 49
$ws.Range("L2").Value2 = 704
$ws.Range("M2").Value2 = 468
$ws.Range("Q2").Value2 = 1125
$ws.Range("C3").Value2 = 405
$ws.Range("E3").Value2 = 59
$ws.Range("F3").Value2 = 237
$ws.Range("G3").Value2 = 95
$ws.Range("H3").Value2 = 59
$ws.Range("I3").Value2 = 140
$ws.Range("J3").Value2 = 144
$ws.Range("N3").Value2 = 23

$ws = $wb.Worksheets.Item("DEF")
$ws.Range("B2").Value2 = 11
$ws.Range("C2").Value2 = 407
$ws.Range("D2").Value2 = 28
$ws.Range("E2").Value2 = 25
$ws.Range("F2").Value2 = 116
$ws.Range("G2").Value2 = 91
$ws.Range("H2").Value2 = 11
$ws.Range("J2").Value2 = 57
$ws.Range("L2").Value2 = 548
$ws.Range("M2").Value2 = 343
$ws.Range("O2").Value2 = 38
$ws.Range("P2").Value2 = 15
$ws.Range("Q2").Value2 = 1052
$ws.Range("C3").Value2 = 292
$ws.Range("D3").Value2 = 6
$ws.Range("E3").Value2 = 74
$ws.Range("F3").Value2 = 199
$ws.Range("G3").Value2 = 60
$ws.Range("H3").Value2 = 73
$ws.Range("I3").Value2 = 110
$ws.Range("J3").Value2 = 96
$ws.Range("N3").Value2 = 64

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value2 = 164
$ws.Range("D2").Value2 = 141
$ws.Range("F2").Value2 = 241
$ws.Range("G2").Value2 = 227
$ws.Range("J2").Value2 = 104
$ws.Range("K2").Value2 = 99
$ws.Range("L2").Value2 = 73
$ws.Range("M2").Value2 = 59
$ws.Range("N2").Value2 = 20
$ws.Range("B3").Value2 = 84

$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("B2").Value2 = 13
$ws.Range("C2").Value2 = 13
$ws.Range("D2").Value2 = 14
$ws.Range("E2").Value2 = 16
$ws.Range("D3").Value2 = 13

$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value2 = 29
$ws.Range("D2").Value2 = 17

